#
# Adds a "Portfolio" link (pointing at the author's new Heroku-hosted
# portfolio site) to the contact-links line at the top of the resume,
# right before the existing "Blog" link, separated by " |  ".
#
# The simplest reliable way to splice new, fully-formed content (a new
# hyperlink field + a plain run) into the middle of an existing
# paragraph in this runtime is to rebuild that paragraph's OOXML in
# full and push it in via Range.InsertXML, which replaces the targeted
# paragraph range with the supplied OOXML and interns any referenced
# relationship ids (creating fresh ones as needed). We keep the
# existing Github / LinkedIn / Blog hyperlinks and their formatting
# byte-for-byte identical, and insert the new "Portfolio" hyperlink
# plus its separator run ahead of "Blog".
#

$d = $word.ActiveDocument

$linksParagraph = $d.Paragraphs.Item(3)
$targetRange = $linksParagraph.Range.Duplicate

$newParagraphXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="00000003"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:hyperlink r:id="rIdGithubKeep"><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:color w:val="1155cc"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Github</w:t></w:r></w:hyperlink><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> | </w:t></w:r><w:hyperlink r:id="rIdLinkedInKeep"><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Roboto" w:cs="Roboto" w:eastAsia="Roboto" w:hAnsi="Roboto"/><w:color w:val="1155cc"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:highlight w:val="white"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">LinkedIn</w:t></w:r></w:hyperlink><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> | </w:t></w:r><w:hyperlink r:id="rIdPortfolioNew"><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:color w:val="1155cc"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Portfolio</w:t></w:r></w:hyperlink><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> |  </w:t></w:r><w:hyperlink r:id="rIdBlogKeep"><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:color w:val="1155cc"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Blog</w:t></w:r></w:hyperlink><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rIdGithubKeep" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://github.com/miodice3" TargetMode="External"/>
<Relationship Id="rIdLinkedInKeep" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.linkedin.com/in/michael-iodice-0158" TargetMode="External"/>
<Relationship Id="rIdBlogKeep" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://dev.to/miodice3" TargetMode="External"/>
<Relationship Id="rIdPortfolioNew" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://miodice3.herokuapp.com/" TargetMode="External"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$targetRange.InsertXML($newParagraphXml)
